# Update loading_percent values for the 380 kV case (rows 2-25, columns B,D,E,F,G,I,K,M)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @{ Row = 2; "B" = 4.237178492551396; "D" = 7.416578437565548; "E" = 10.14009417437717; "F" = 40.64615495779413; "G" = 3.685806719290623; "I" = 25.61414345103336; "K" = 14.28495389079682; "M" = 16.0924152573312 }
    @{ Row = 3; "B" = 4.168288532743712; "D" = 7.400255215415634; "E" = 9.984895581274426; "F" = 40.03323973719511; "G" = 3.689995067957691; "I" = 25.47230243231675; "K" = 13.98794304434639; "M" = 15.89257095372307 }
    @{ Row = 4; "B" = 4.124031337052815; "D" = 7.390623555279588; "E" = 9.891372203056402; "F" = 39.66115592446533; "G" = 3.692696389444505; "I" = 25.38702023837573; "K" = 13.80812069075922; "M" = 15.77409617897681 }
    @{ Row = 5; "B" = 4.105510053816168; "D" = 7.386796758811228; "E" = 9.853750431163027; "F" = 39.51076047170724; "G" = 3.693829940469093; "I" = 25.35273403173654; "K" = 13.73560085840645; "M" = 15.72693719770711 }
    @{ Row = 6; "B" = 4.102405442583064; "D" = 7.386167247791141; "E" = 9.847534236354328; "F" = 39.48586654257959; "G" = 3.694020147218499; "I" = 25.3470693762988; "K" = 13.72360835137514; "M" = 15.71917576054591 }
    @{ Row = 7; "B" = 4.123783512151402; "D" = 7.39057154843237; "E" = 9.890862781063463; "F" = 39.65912243775817; "G" = 3.69271154416592; "I" = 25.38655593651728; "K" = 13.80713943064541; "M" = 15.77345556696256 }
    @{ Row = 8; "B" = 4.213835399396049; "D" = 7.410869021174891; "E" = 10.08624172049441; "F" = 40.43403859555357; "G" = 3.68722403393939; "I" = 25.56486877769344; "K" = 14.182088262252; "M" = 16.02266570239525 }
    @{ Row = 9; "B" = 4.374630398427191; "D" = 7.4537772090098; "E" = 10.48140748261302; "F" = 41.98011511419422; "G" = 3.677485570255045; "I" = 25.92832844622258; "K" = 14.93233898854605; "M" = 16.542325813829 }
    @{ Row = 10; "B" = 4.482823579921749; "D" = 7.487198259808451; "E" = 10.77635915015993; "F" = 43.12246824018072; "G" = 3.670945312421828; "I" = 26.20297544750568; "K" = 15.48588394331281; "M" = 16.93938246305247 }
    @{ Row = 11; "B" = 4.529826173890804; "D" = 7.502813616822031; "E" = 10.91097476151062; "F" = 43.64153608651222; "G" = 3.668101552080046; "I" = 26.32939008141807; "K" = 15.73685197295095; "M" = 17.1225415553793 }
    @{ Row = 12; "B" = 4.54730186266878; "D" = 7.508785745978829; "E" = 10.96196709471085; "F" = 43.8378440009546; "G" = 3.667043448837222; "I" = 26.37745570515941; "K" = 15.83165958490848; "M" = 17.19219791804348 }
    @{ Row = 13; "B" = 4.543552603887667; "D" = 7.507496925687778; "E" = 10.95098499225233; "F" = 43.79557937354983; "G" = 3.667270497750483; "I" = 26.36709551237139; "K" = 15.8112528790794; "M" = 17.17718397522838 }
    @{ Row = 14; "B" = 4.531270411577198; "D" = 7.503303768391241; "E" = 10.91516986904433; "F" = 43.65769243848472; "G" = 3.668014125932573; "I" = 26.33334065825652; "K" = 15.74465702609883; "M" = 17.12826663364729 }
    @{ Row = 15; "B" = 4.523704995466293; "D" = 7.500743002928732; "E" = 10.8932328239831; "F" = 43.57319506383462; "G" = 3.668472059740798; "I" = 26.31268973140165; "K" = 15.70383229669987; "M" = 17.09834018863449 }
    @{ Row = 16; "B" = 4.479706871558678; "D" = 7.486186024495886; "E" = 10.76756651987789; "F" = 43.0885196451842; "G" = 3.671133792279203; "I" = 26.19474220745777; "K" = 15.46945632451265; "M" = 16.9274579026789 }
    @{ Row = 17; "B" = 4.452144697274973; "D" = 7.477361105415559; "E" = 10.69055330669858; "F" = 42.79091488470991; "G" = 3.672800247982212; "I" = 26.12275131349507; "K" = 15.3253802150002; "M" = 16.82322907339194 }
    @{ Row = 18; "B" = 4.436083429353372; "D" = 7.472324043214813; "E" = 10.64630156089595; "F" = 42.61969423322713; "G" = 3.673771128538768; "I" = 26.08148407692434; "K" = 15.2424382252398; "M" = 16.7635220813009 }
    @{ Row = 19; "B" = 4.430609743421856; "D" = 7.470625252237014; "E" = 10.6313277282481; "F" = 42.56171906401848; "G" = 3.674101982231339; "I" = 26.06753626237347; "K" = 15.21434628664976; "M" = 16.74335004037717 }
    @{ Row = 20; "B" = 4.455100324112136; "D" = 7.478296519683127; "E" = 10.69874726317606; "F" = 42.82260142834065; "G" = 3.672621570681729; "I" = 26.13040049793779; "K" = 15.34072573812387; "M" = 16.83429977944692 }
    @{ Row = 21; "B" = 4.534886795313654; "D" = 7.504533803709509; "E" = 10.92568956422875; "F" = 43.69820130869124; "G" = 3.667795196037011; "I" = 26.34325012014417; "K" = 15.76422487072891; "M" = 17.14262728336261 }
    @{ Row = 22; "B" = 4.585146258233221; "D" = 7.522024697654094; "E" = 11.07408379059963; "F" = 44.26891524455863; "G" = 3.664750208383965; "I" = 26.48348880973385; "K" = 16.03962850103606; "M" = 17.3458436940227 }
    @{ Row = 23; "B" = 4.558495814976717; "D" = 7.512658184306531; "E" = 10.99489138496955; "F" = 43.96450911536405; "G" = 3.666365416132831; "I" = 26.40854325781342; "K" = 15.892800345391; "M" = 17.23724895223557 }
    @{ Row = 24; "B" = 4.453764755806062; "D" = 7.477873505176742; "E" = 10.69504269656576; "F" = 42.80827631462219; "G" = 3.672702310769817; "I" = 26.12694192257402; "K" = 15.33378836375698; "M" = 16.8292940367828 }
    @{ Row = 25; "B" = 4.332856070333254; "D" = 7.441836348881807; "E" = 10.37349711007585; "F" = 41.56003177063429; "G" = 3.680011521375251; "I" = 25.82860577686017; "K" = 14.72848804288949; "M" = 16.39879780703253 }
)

foreach ($entry in $newValues) {
    foreach ($col in @('B', 'D', 'E', 'F', 'G', 'I', 'K', 'M')) {
        $addr = "$col$($entry.Row)"
        $ws.Range($addr).Value = $entry[$col]
    }
}

Write-Host "Updated $(($newValues | Measure-Object).Count) rows across columns B,D,E,F,G,I,K,M"
